# Update cryptos list cell values (Price + Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.825.47'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.634.93'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'215.00"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('E10').Value = '  +2.94%  '
$ws.Range('D11').Value = "'0.0780"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').Value = '1.638.07'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').Value = '1.860.04'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('D17').Value = "'63.05"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '25.827.53'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').Value = "'193.85"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('E21').Value = '  +1.70%  '
$ws.Range('D22').Value = "'9.93"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.22%  '
$ws.Range('D23').Value = "'6.17"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.99%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('E25').Value = '  -1.90%  '
$ws.Range('D26').Value = "'139.34"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.03%  '
$ws.Range('E27').Value = '  -3.78%  '
$ws.Range('E28').Value = '  +1.60%  '
$ws.Range('D29').Value = "'15.46"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.50%  '
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('D31').Value = "'0.0495"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.73%  '
$ws.Range('E32').Value = '  +1.42%  '
$ws.Range('E33').Value = '  +1.73%  '
$ws.Range('E34').Value = '  +1.45%  '
$ws.Range('E35').Value = '  +0.41%  '
$ws.Range('D36').Value = "'0.902"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.92%  '
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('D38').Value = "'0.551"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.76%  '
$ws.Range('D39').Value = '1.116.77'
$ws.Range('E39').Value = '  -1.01%  '
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('D43').Value = "'99.36"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.14%  '
$ws.Range('E44').Value = '  +0.34%  '
$ws.Range('E45').Value = '  -2.83%  '
$ws.Range('D46').Value = "'55.53"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.68%  '
$ws.Range('E47').Value = '  +11.07%  '
$ws.Range('E48').Value = '  -5.26%  '
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('D50').Value = "'7.61"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('E51').Value = '  +0.20%  '
